$d = $word.ActiveDocument

$findRange = $d.Content
$found = $findRange.Find.Execute(" per dire riesegui il comando una seconda volta mettendo nel secondo file l’output con l’input del primo) e lo si può far seguire da “tee”.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    Write-Output "ERROR: original text not found"
} else {
    $insPoint = $findRange.Start
    $delEnd = $findRange.End

    # Delete the old run's text entirely, leaving a collapsed insertion point
    $findRange.Text = ""
    $cursor = $insPoint

    $r1 = $d.Range($cursor, $cursor)
    $r1.InsertAfter(" per dire esegui il ")
    $b1 = $r1.End

    $r2 = $d.Range($b1, $b1)
    $r2.InsertAfter("primo ")
    $b2 = $r2.End

    $r3 = $d.Range($b2, $b2)
    $r3.InsertAfter("comando mettendo nel secondo l’output del primo")
    $b3 = $r3.End

    $r4 = $d.Range($b3, $b3)
    $r4.InsertAfter(", tuttavia avviali assieme")
    $b4 = $r4.End

    $r5 = $d.Range($b4, $b4)
    $r5.InsertAfter(") e lo si può far seguire da “tee”.")
    $b5 = $r5.End

    # Force run separation between the 5 newly inserted pieces (and from the
    # untouched run that precedes them) by toggling a character property
    # across each full span after all text is in place. (Toggling and
    # reverting on the whole span keeps the run boundary without leaving a
    # visible formatting change.)
    $t1 = $d.Range($insPoint, $b1)
    $t1.Bold = 1
    $t1.Bold = 0

    $t2 = $d.Range($b1, $b2)
    $t2.Bold = 1
    $t2.Bold = 0

    $t3 = $d.Range($b2, $b3)
    $t3.Bold = 1
    $t3.Bold = 0

    $t4 = $d.Range($b3, $b4)
    $t4.Bold = 1
    $t4.Bold = 0

    $t5 = $d.Range($b4, $b5)
    $t5.Bold = 1
    $t5.Bold = 0

    Write-Output "Replacement complete"
}
